# Daily attendance processing - swap the order of the two comma-separated
# "Recorded By" identities in column G for the specific rows touched by
# this run (System <-> user email, or admin <-> user email), leaving any
# single-value or three-value (e.g. backup@backdoor.com, System, system)
# cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,6,7,10,12,13,14,15,18,19,20,21,22,24,26,29,32,33,36,38,39,40,41,44,45,46,47,48,50,52,55,58,59,62,64,65,66,67,70,71,72,73,74,76,78,83,84,85,86,87,90,92,99,101,109,110,111,112,113,116,118,125,127,135,136,137,138,139,142,144,151,153)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $value = $cell.Value2
    $parts = $value -split ', ', 2
    if ($parts.Length -eq 2) {
        $cell.Value2 = $parts[1] + ', ' + $parts[0]
    }
}
